$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{Row=5; I='sv'; J='Statement-opinion'}
    @{Row=7; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=8; I='sd'; J='Statement-non-opinion'}
    @{Row=29; I='sv'; J='Statement-opinion'}
    @{Row=46; I='sd'; J='Statement-non-opinion'}
    @{Row=56; I='sd'; J='Statement-non-opinion'}
    @{Row=58; I='ba'; J='Appreciation'}
    @{Row=67; I='sd'; J='Statement-non-opinion'}
    @{Row=72; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=85; I='sd'; J='Statement-non-opinion'}
    @{Row=89; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=90; I='sv'; J='Statement-opinion'}
    @{Row=92; I='sv'; J='Statement-opinion'}
    @{Row=96; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=100; I='sv'; J='Statement-opinion'}
    @{Row=104; I='sd'; J='Statement-non-opinion'}
    @{Row=109; I='sv'; J='Statement-opinion'}
    @{Row=114; I='sv'; J='Statement-opinion'}
    @{Row=165; I='aa'; J='Agree/Accept'}
    @{Row=176; I='sd'; J='Statement-non-opinion'}
    @{Row=179; I='sd'; J='Statement-non-opinion'}
    @{Row=201; I='sv'; J='Statement-opinion'}
    @{Row=212; I='sd'; J='Statement-non-opinion'}
    @{Row=251; I='sv'; J='Statement-opinion'}
    @{Row=252; I='aa'; J='Agree/Accept'}
    @{Row=264; I='sd'; J='Statement-non-opinion'}
    @{Row=270; I='sv'; J='Statement-opinion'}
    @{Row=284; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=291; I='sv'; J='Statement-opinion'}
    @{Row=296; I='aa'; J='Agree/Accept'}
    @{Row=297; I='sv'; J='Statement-opinion'}
    @{Row=308; I='ba'; J='Appreciation'}
    @{Row=332; I='sd'; J='Statement-non-opinion'}
    @{Row=340; I='sd'; J='Statement-non-opinion'}
    @{Row=343; I='aa'; J='Agree/Accept'}
    @{Row=345; I='%'; J='Uninterpretable'}
    @{Row=361; I='sv'; J='Statement-opinion'}
    @{Row=377; I='sv'; J='Statement-opinion'}
    @{Row=383; I='sv'; J='Statement-opinion'}
    @{Row=387; I='sv'; J='Statement-opinion'}
    @{Row=391; I='ba'; J='Appreciation'}
    @{Row=413; I='sv'; J='Statement-opinion'}
    @{Row=416; I='sd'; J='Statement-non-opinion'}
    @{Row=422; I='aa'; J='Agree/Accept'}
    @{Row=424; I='sd'; J='Statement-non-opinion'}
    @{Row=437; I='sv'; J='Statement-opinion'}
    @{Row=458; I='qy'; J='Yes-No-Question'}
    @{Row=471; I='aa'; J='Agree/Accept'}
    @{Row=481; I='sv'; J='Statement-opinion'}
    @{Row=482; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=495; I='sd'; J='Statement-non-opinion'}
    @{Row=508; I='aa'; J='Agree/Accept'}
    @{Row=513; I='sd'; J='Statement-non-opinion'}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}

Write-Output ("Updated {0} rows" -f $updates.Count)